$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91828
$ws.Range("B3").Value = 99013
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 79267
